$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.618.45'
$ws.Range("E2").Value = '  -1.25%  '

$ws.Range("D3").Value = '3.678.86'
$ws.Range("E3").Value = '  -3.04%  '

$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.97'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +2.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.27'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +11.84%  '

$ws.Range("D7").Value = '3.678.32'
$ws.Range("E7").Value = '  -3.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.629'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -3.66%  '

$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.716'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  -1.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.162'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -4.99%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.91'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +8.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000291'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  -7.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.62'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  -3.44%  '

$ws.Range("D15").Value = '4.287.48'
$ws.Range("E15").Value = '  -2.22%  '

$ws.Range("D16").Value = '3.702.52'
$ws.Range("E16").Value = '  -2.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.30'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -4.99%  '

$ws.Range("E18").Value = '  -1.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.80'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  -5.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.12'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  -5.88%  '

$ws.Range("D21").Value = '68.516.44'
$ws.Range("E21").Value = '  -1.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '409.91'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  -4.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.57'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -1.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.70'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -3.25%  '

$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.10'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +5.24%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.03'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -5.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.77'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -5.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.88'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -1.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.05'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +1.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.44'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -7.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.65'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  -4.61%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.25'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -9.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.37'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -6.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.118'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  -3.67%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '65.05'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  -4.54%  '

$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '607.04'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -3.15%  '

$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '43.41'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -7.35%  '

$ws.Range("D38").Value = '0.0₃0894'
$ws.Range("E38").Value = '  -7.91%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.400'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -4.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -0.10%  '

$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.136'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -4.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.03'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -4.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.72'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -4.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.94'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -8.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0438'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  -4.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.12'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -5.53%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.20'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.94%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.134'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  -4.53%  '

$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.71'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -1.29%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.757.58'
$ws.Range("E51").Value = '  -2.13%  '
